$d = $word.ActiveDocument

$replacements = @(
    @("62×48=2976", "24×79=1896"),
    @("87×34=2958", "70×11=770"),
    @("66×20=1320", "26×99=2574"),
    @("76×38=2888", "11×47=517"),
    @("51×62=3162", "89×56=4984"),
    @("33×26=858",  "58×68=3944"),
    @("13×83=1079", "79×89=7031"),
    @("38×90=3420", "99×40=3960"),
    @("53×43=2279", "67×32=2144"),
    @("37×18=666",  "33×67=2211"),
    @("41×83=3403", "67×82=5494"),
    @("55×82=4510", "45×66=2970"),
    @("98×65=6370", "95×94=8930"),
    @("47×84=3948", "86×25=2150"),
    @("91×30=2730", "71×88=6248"),
    @("68×95=6460", "59×39=2301"),
    @("17×88=1496", "92×94=8648"),
    @("72×40=2880", "35×63=2205"),
    @("91×66=6006", "45×62=2790"),
    @("64×87=5568", "43×14=602"),
    @("51×73=3723", "50×39=1950"),
    @("11×72=792",  "53×27=1431"),
    @("44×68=2992", "40×19=760"),
    @("68×31=2108", "54×76=4104"),
    @("46×51=2346", "45×41=1845")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
